# API_END_POINTS.xlsx - "Added Changes for API ShortVisit"
# Rebuilds the endpoints table with new /:shorty and geo_detail pagination rows,
# re-orders existing rows, and re-applies the (re-themed) row colouring.

function RGBColor($r, $g, $b) {
    # Excel COM colors are encoded 0x00BBGGRR
    return ($b * 65536) + ($g * 256) + $r
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Table content (A1:D11)
# ---------------------------------------------------------------------------
# Note: B7 ("GET[HTTP/API]") is intentionally left blank here and written
# afterwards, so that it becomes the last *new* string appended to
# sharedStrings.xml (matching the order in which the table was originally
# edited), after /:shorty, /:shorty/geo_detail?page=:page_num, and the
# accompanying ":shorty => ... :page_num => ..." description.
$rows = @(
    @("API URL", "METHOD", "DATA", "Response"),
    @("/users", "POST", "{:email,:password,:name}", "User Name with API TOKEN"),
    @("/short_urls", "POST", "{:original_url}", "Short URL Details"),
    @("/users/regenerate_token", "POST", "{:email,:password,:name}", "NEW APITOKEN"),
    @("/short_urls/:shorty", "DELETE", ":shorty => Short URL", "Deletion Status"),
    @("/short_urls?page=1", "GET", "N/A", "LIST OF SHORTURLS Paginated : Perpage 10"),
    @("/:shorty", $null, ":shorty => Short URL", "Displays the Original URL"),
    @("/short_urls", "GET", "N/A", "LIST OF SHORTURLS"),
    @("/:shorty/geo_detail", "GET", ":shorty => Short URL", "GEO Location Details"),
    @("/:shorty/geo_detail?page=:page_num", "GET", ":shorty => Short URL, :page_num => Current page Number", "GEO Location Details"),
    @("/users/authenticate", "GET", "Authorization: Token APITOKEN", "Authorized")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 1
    $rowVals = $rows[$i]
    $ws.Range("A$r").Value2 = $rowVals[0]
    if ($null -ne $rowVals[1]) {
        $ws.Range("B$r").Value2 = $rowVals[1]
    }
    $ws.Range("C$r").Value2 = $rowVals[2]
    $ws.Range("D$r").Value2 = $rowVals[3]
}

$ws.Range("B7").Value2 = "GET[HTTP/API]"

# ---------------------------------------------------------------------------
# 2. Formatting
# ---------------------------------------------------------------------------
$darkGray = RGBColor 0x33 0x33 0x33
$yellow   = RGBColor 0xFF 0xFF 0x00
$lightGreen = RGBColor 0xCC 0xFF 0xCC

# A11 used to be styled off the built-in "Hyperlink" cell style (underlined);
# clear that leftover formatting before re-applying the new row styling below.
$ws.Range("A11").Font.Underline = -4142

# Row 1 - header: bold 16pt Helvetica yellow text on accent1 (theme4) fill
$hdr = $ws.Range("A1:D1")
$hdr.Font.Name = "Helvetica"
$hdr.Font.Size = 16
$hdr.Font.Bold = $true
$hdr.Font.Underline = -4142
$hdr.Font.Color = $yellow
$hdr.Interior.ThemeColor = 5
$hdr.Interior.Pattern = 1

# Rows 2-4 - 16pt Helvetica dark-gray text on light accent1 (theme4) fill
$body1 = $ws.Range("A2:D4")
$body1.Font.Name = "Helvetica"
$body1.Font.Size = 16
$body1.Font.Bold = $false
$body1.Font.Underline = -4142
$body1.Font.Color = $darkGray
$body1.Interior.ThemeColor = 5
$body1.Interior.Pattern = 1

# Row 5 - 16pt Helvetica dark-gray text on accent6 (theme9) fill
$body2 = $ws.Range("A5:D5")
$body2.Font.Name = "Helvetica"
$body2.Font.Size = 16
$body2.Font.Bold = $false
$body2.Font.Underline = -4142
$body2.Font.Color = $darkGray
$body2.Interior.ThemeColor = 10
$body2.Interior.Pattern = 1

# Rows 6-11 - 16pt Helvetica dark-gray text on light-green fill
$body3 = $ws.Range("A6:D11")
$body3.Font.Name = "Helvetica"
$body3.Font.Size = 16
$body3.Font.Bold = $false
$body3.Font.Underline = -4142
$body3.Font.Color = $darkGray
$body3.Interior.Color = $lightGreen

# D10 keeps the light-green fill, but reverts to the plain default font
$d10 = $ws.Range("D10")
$d10.Font.Name = "Calibri"
$d10.Font.Size = 12
$d10.Font.Bold = $false
$d10.Font.Underline = -4142
$d10.Font.ThemeColor = 2
$d10.Interior.Color = $lightGreen

# Row 11 is a brand-new data row (the table used to stop at row 10); give it
# the same row height as the rest of the table.
$ws.Rows.Item(11).RowHeight = 17

# ---------------------------------------------------------------------------
# 3. View / print setup
# ---------------------------------------------------------------------------
$ws.Range("C19").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
